$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New "Participant ID" query text (cell B2) - replaces the old query that
# is no longer referenced anywhere, so the old shared string drops out and
# this new text gets appended to the shared-string table.
$participantQuery = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE g.instrument_model in ['Illumina MiSeq']
WITH p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id limit 100
'@

$ws.Cells.Item(2, 2).Value = $participantQuery

# Save so the current selection / view settings match the target state.
$ws.Range("C3").Select()
